# Adds the four missing "t / tt" glyph textboxes to slide 6 of the TeeTime
# logo deck (mirrors the existing "T" + "=" glyph pair already on the slide).
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points (Single,
# i.e. 32-bit float) while OOXML stores EMU (1 pt = 12700 EMU). A naive
# emu/12700.0 assignment loses precision once the Single truncates, so
# EmuToPtExact searches (in the same Single precision the host uses) for a
# points value whose float32 round-trip reproduces the exact target EMU.
function EmuToPtExact($emuTarget) {
    $pts = [double]$emuTarget / 12700.0
    $step = 0.0000005
    $guard = 0
    while ($guard -lt 200000) {
        $f32 = [single]$pts
        $back = [int]([double]$f32 * 12700.0)
        if ($back -eq $emuTarget) {
            return $pts
        }
        if ($back -lt $emuTarget) {
            $pts = $pts + $step
        } else {
            $pts = $pts - $step
        }
        $guard = $guard + 1
    }
    return $pts
}

function SetShapePosition($shape, $xEmu, $yEmu, $cxEmu, $cyEmu) {
    $shape.Left = EmuToPtExact $xEmu
    $shape.Top = EmuToPtExact $yEmu
    $shape.Width = EmuToPtExact $cxEmu
    $shape.Height = EmuToPtExact $cyEmu
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# Existing shapes used as style donors for the new ones:
#   Shapes.Item(1) = "TextBox 1" ("eeTime") -> wrap="none" glyph style (t / tt)
#   Shapes.Item(3) = "TextBox 3" ("=")      -> wrap="square" glyph style (=)
$wordTemplate = $s.Shapes.Item(1)
$barTemplate = $s.Shapes.Item(3)

# TextBox 4 (id 5): "=" glyph
$dup = $barTemplate.Duplicate()
$shp = $dup.Item(1)
$shp.Name = "TextBox 4"
SetShapePosition $shp 8417808 4557500 413578 707886
$shp.TextFrame.TextRange.Text = "≡"

# TextBox 5 (id 6): "t" glyph
$dup = $wordTemplate.Duplicate()
$shp = $dup.Item(1)
$shp.Name = "TextBox 5"
SetShapePosition $shp 8230446 4556768 263855 615553
$shp.TextFrame.TextRange.Text = "t"

# TextBox 6 (id 7): "=" glyph
$dup = $barTemplate.Duplicate()
$shp = $dup.Item(1)
$shp.Name = "TextBox 6"
SetShapePosition $shp 6866454 5173053 413578 707886
$shp.TextFrame.TextRange.Text = "≡"

# TextBox 7 (id 8): "tt" glyph
$dup = $wordTemplate.Duplicate()
$shp = $dup.Item(1)
$shp.Name = "TextBox 7"
SetShapePosition $shp 6521554 5173053 428131 615553
$shp.TextFrame.TextRange.Text = "tt"
